$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70: Consecrating Congregation / Holy Water
$ws.Range("H70").Value = 3127.85
$ws.Range("J70").Value = 3761.3076
$ws.Range("L70").Value = 11283.9228
$ws.Range("N70").Value = -11823.9228

# Row 73: Curbing the Contagion (L) / Holy Water
$ws.Range("H73").Value = 3127.85
$ws.Range("J73").Value = 3761.3076
$ws.Range("L73").Value = 11283.9228
$ws.Range("N73").Value = -13155.9228

# Row 88: The Grave of Hemlock Groves / Growth Formula Zeta
$ws.Range("H88").Value = 501.8
$ws.Range("I88").Value = 449.66666
$ws.Range("J88").Value = 524.1429000000001
$ws.Range("K88").Value = 449.66666
$ws.Range("L88").Value = 524.1429000000001
$ws.Range("M88").Value = -43.66665999999998
$ws.Range("N88").Value = -1336.1429

# Row 91: Dappling the Highlands (L) / Growth Formula Zeta
$ws.Range("H91").Value = 501.8
$ws.Range("I91").Value = 449.66666
$ws.Range("J91").Value = 524.1429000000001
$ws.Range("K91").Value = 449.66666
$ws.Range("L91").Value = 524.1429000000001
$ws.Range("M91").Value = 954.33334
$ws.Range("N91").Value = -3332.1429

# Row 100: Asking for a Friend / Beetle Glue
$ws.Range("H100").Value = 1669.25
$ws.Range("I100").Value = 1669.25
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1669.25
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -1128.25

# Row 111: An Eye for Healing / Grade 1 Dexterity Alkahest
$ws.Range("H111").Value = 1850
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 1850
$ws.Range("K111").Value = 0
$ws.Range("L111").ClearContents()
$ws.Range("M111").Value = 5550
$ws.Range("N111").Value = -11684

# Row 125: Body over Mind / Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 3383.4285
$ws.Range("I125").Value = 2976.8
$ws.Range("J125").Value = 4400
$ws.Range("K125").Value = 26791.2
$ws.Range("L125").Value = 39600
$ws.Range("M125").Value = -24331.2
$ws.Range("N125").Value = -44520

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 4694.696
$ws.Range("I137").Value = 1997.9
$ws.Range("K137").Value = 5993.700000000001
$ws.Range("M137").Value = -3443.700000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 8: You've Got Mail / Bronze Haubergeon
$ws.Range("H8").Value = 3184
$ws.Range("J8").Value = 3783
$ws.Range("L8").Value = 3783
$ws.Range("N8").Value = -4071

# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 1968
$ws.Range("I45").Value = 1349
$ws.Range("K45").Value = 1349
$ws.Range("M45").Value = -972

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 4596.2607
$ws.Range("I61").Value = 3435
$ws.Range("J61").Value = 5863.091
$ws.Range("K61").Value = 3435
$ws.Range("L61").Value = 5863.091
$ws.Range("M61").Value = -3223
$ws.Range("N61").Value = -6287.091

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 1814.3182
$ws.Range("I74").Value = 1424.5238
$ws.Range("K74").Value = 1424.5238
$ws.Range("M74").Value = -550.5237999999999

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 1814.3182
$ws.Range("I77").Value = 1424.5238
$ws.Range("K77").Value = 7122.619
$ws.Range("M77").Value = -2754.619

# Row 125: The Incomplete Costume / High Durium Armor of Fending
$ws.Range("H125").Value = 31750.25
$ws.Range("J125").Value = 31750.25
$ws.Range("L125").Value = 31750.25
$ws.Range("N125").Value = -41590.25

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2856.8215
$ws.Range("I132").Value = 2444
$ws.Range("K132").Value = 7332
$ws.Range("M132").Value = -4802

# Row 135: Forgiveness for My Shins / Ruthenium Sabatons of Fending
$ws.Range("H135").Value = 79999
$ws.Range("J135").Value = 79999
$ws.Range("L135").Value = 79999
$ws.Range("N135").Value = -90139

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 4596.2607
$ws.Range("I136").Value = 3435
$ws.Range("J136").Value = 5863.091
$ws.Range("K136").Value = 10305
$ws.Range("L136").Value = 17589.273
$ws.Range("M136").Value = -7755
$ws.Range("N136").Value = -22689.273

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 3172
$ws.Range("I20").Value = 1888
$ws.Range("K20").Value = 1888
$ws.Range("M20").Value = -1641

# Row 25: Tools of the Trade / Iron Doming Hammer
$ws.Range("H25").Value = 2103.5
$ws.Range("I25").Value = 3007
$ws.Range("J25").Value = 1200
$ws.Range("K25").Value = 3007
$ws.Range("L25").Value = 1200
$ws.Range("M25").Value = -2772
$ws.Range("N25").Value = -1670

# Row 37: That's Some Fine Grinding / Initiate's Mortar
$ws.Range("H37").Value = 1761
$ws.Range("I37").Value = 1307.7142
$ws.Range("J37").Value = 2214.2856
$ws.Range("K37").Value = 1307.7142
$ws.Range("L37").Value = 2214.2856
$ws.Range("M37").Value = -1170.7142
$ws.Range("N37").Value = -2488.2856

# Row 39: Out on a Limb / Bas-relief Steel Saw
$ws.Range("H39").Value = 19833
$ws.Range("J39").Value = 19833
$ws.Range("L39").Value = 19833
$ws.Range("N39").Value = -20611

# Row 54: Get Me to the War on Time / Cobalt Doming Hammer
$ws.Range("H54").Value = 8617.200000000001
$ws.Range("I54").Value = 2021.5
$ws.Range("J54").Value = 35000
$ws.Range("K54").Value = 2021.5
$ws.Range("L54").Value = 35000
$ws.Range("M54").Value = -1537.5
$ws.Range("N54").Value = -35968

# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 2025
$ws.Range("I99").Value = 2025
$ws.Range("K99").Value = 2025
$ws.Range("M99").Value = -527

$ws = $wb.Worksheets.Item("CRP")
# Row 15: On the Move / Ragstone Grinding Wheel
$ws.Range("H15").Value = 4933
$ws.Range("I15").Value = 104
$ws.Range("J15").Value = 7347.5
$ws.Range("K15").Value = 104
$ws.Range("L15").Value = 7347.5
$ws.Range("M15").Value = 66
$ws.Range("N15").Value = -7687.5

# Row 50: The Arsenal of Theocracy / Cobalt Halberd
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("N50").Value = 0

# Row 57: Clogs of War / Mahogany Pattens
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

# Row 98: Pinewheel / Pine Spinning Wheel
$ws.Range("H98").Value = 40499.75
$ws.Range("J98").Value = 40499.75
$ws.Range("L98").Value = 40499.75
$ws.Range("N98").Value = -44991.75

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 1851.7693
$ws.Range("I134").Value = 1851.7693
$ws.Range("K134").Value = 5555.3079
$ws.Range("M134").Value = -3020.3079

$ws = $wb.Worksheets.Item("GSM")
# Row 32: Love in the Time of Umbra / Silver Ring
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").ClearContents()
$ws.Range("N32").Value = 0

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 1404.16
$ws.Range("I102").Value = 1141.1364
$ws.Range("J102").Value = 3333
$ws.Range("K102").Value = 1141.1364
$ws.Range("L102").Value = 3333
$ws.Range("M102").Value = 480.8635999999999
$ws.Range("N102").Value = -6577

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 1804.6666
$ws.Range("I122").Value = 1534.9
$ws.Range("J122").Value = 2344.2
$ws.Range("K122").Value = 4604.700000000001
$ws.Range("L122").Value = 7032.599999999999
$ws.Range("M122").Value = -2154.700000000001
$ws.Range("N122").Value = -11932.6

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

$ws = $wb.Worksheets.Item("LTW")
# Row 21: Heads Up / Hard Leather Skullcap
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").ClearContents()
$ws.Range("N21").Value = 0

# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 620
$ws.Range("I22").Value = 597
$ws.Range("J22").Value = 643
$ws.Range("K22").Value = 597
$ws.Range("L22").Value = 643
$ws.Range("M22").Value = -302
$ws.Range("N22").Value = -1233

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 620
$ws.Range("I27").Value = 597
$ws.Range("J27").Value = 643
$ws.Range("K27").Value = 597
$ws.Range("L27").Value = 643
$ws.Range("M27").Value = -490
$ws.Range("N27").Value = -857

# Row 100: Tiger in the Sack / Tiger Leather
$ws.Range("H100").Value = 5996.357
$ws.Range("I100").Value = 3772.111
$ws.Range("K100").Value = 3772.111
$ws.Range("M100").Value = -3231.111

$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display / Ruby Cotton Cloth
$ws.Range("H96").Value = 3003
$ws.Range("I96").Value = 3003
$ws.Range("K96").Value = 3003
$ws.Range("M96").Value = -1630

# Row 100: Of Great Import / Kudzu Thread
$ws.Range("H100").Value = 624.375
$ws.Range("I100").Value = 624.375
$ws.Range("K100").Value = 1248.75
$ws.Range("M100").Value = -707.75

# Row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 758
$ws.Range("I107").Value = 599.6667
$ws.Range("K107").Value = 1799.0001
$ws.Range("M107").Value = 120.9999

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1848.3889
$ws.Range("I132").Value = 1780.6471
$ws.Range("K132").Value = 5341.9413
$ws.Range("M132").Value = -2811.9413
